# Apply CBO data revisions ("Full CBO changes after talking to Louise")
# Updates the "current" section (rows 6, 12, 16) and the corresponding
# "difference" section (rows 34, 40, 44) on Sheet 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Row 6: Federal Corporate Taxes Contribution / current
$ws.Range("R6").Value = -0.0154

# Row 12: Federal Social Benefits Contribution / current
$ws.Range("F12").Value = -0.9172
$ws.Range("G12").Value = -0.2334
$ws.Range("H12").Value = -0.8409
$ws.Range("I12").Value = -0.4083
$ws.Range("J12").Value = 0.0142
$ws.Range("K12").Value = 0.124
$ws.Range("L12").Value = 0.1722
$ws.Range("M12").Value = -0.1072
$ws.Range("N12").Value = -0.1102
$ws.Range("O12").Value = -0.095
$ws.Range("P12").Value = -0.08
$ws.Range("Q12").Value = -0.0445
$ws.Range("R12").Value = -0.0796

# Row 16: Fiscal Impact / current
$ws.Range("F16").Value = -3.8339
$ws.Range("G16").Value = -4.8217
$ws.Range("H16").Value = -2.4774
$ws.Range("I16").Value = -1.0032
$ws.Range("J16").Value = -1.5351
$ws.Range("K16").Value = -1.6397
$ws.Range("L16").Value = -0.574
$ws.Range("M16").Value = -0.5699
$ws.Range("N16").Value = -0.6183
$ws.Range("O16").Value = -0.1956
$ws.Range("P16").Value = -0.1398
$ws.Range("Q16").Value = 0.0254
$ws.Range("R16").Value = -0.052

# Row 34: Federal Corporate Taxes Contribution / difference
$ws.Range("R34").Value = 0.0312

# Row 40: Federal Social Benefits Contribution / difference
$ws.Range("F40").Value = 0.0092
$ws.Range("G40").Value = 0.0089
$ws.Range("H40").Value = 0.0086
$ws.Range("I40").Value = 0.0205
$ws.Range("J40").Value = -0.0038
$ws.Range("K40").Value = -0.0161
$ws.Range("L40").Value = -0.0274
$ws.Range("M40").Value = -0.0558
$ws.Range("N40").Value = -0.0439
$ws.Range("O40").Value = -0.0312
$ws.Range("P40").Value = -0.0184
$ws.Range("Q40").Value = -0.0036
$ws.Range("R40").Value = -0.0013

# Row 44: Fiscal Impact / difference
$ws.Range("F44").Value = 0.0092
$ws.Range("G44").Value = 0.0089
$ws.Range("H44").Value = 0.0086
$ws.Range("I44").Value = 0.0205
$ws.Range("J44").Value = -0.0632
$ws.Range("K44").Value = 0.1503
$ws.Range("L44").Value = 0.3096
$ws.Range("M44").Value = 0.1331
$ws.Range("N44").Value = 0.1242
$ws.Range("O44").Value = 0.1485
$ws.Range("P44").Value = 0.196
$ws.Range("Q44").Value = 0.1191
$ws.Range("R44").Value = -0.0862
